# Insert a new bullet paragraph "Hoe veel algoritmes moeten we ongeveer
# testen?" right after the "Testen: Hoe moeten de CSV files ingelezen
# worden?" paragraph, and move the (hidden) "_GoBack" bookmark from the
# end of the old last paragraph onto the end of the newly added one.

$d = $word.ActiveDocument

# The "Testen..." paragraph is currently the last paragraph in the body.
$lastPara = $d.Paragraphs.Last

# Add a new paragraph right after it; it inherits the same pPr
# (NoSpacing style + numPr ilvl=0/numId=2) from the paragraph it split off.
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last

# A trailing sentinel character is appended first so the bookmark can be
# anchored safely in the middle of the run text (anchoring a collapsed
# range exactly at the end of a paragraph's text mis-serializes the
# bookmark). The sentinel is stripped afterwards, leaving the bookmark
# collapsed right after the run, matching the original placement.
$newPara.Range.Text = "Hoe veel algoritmes moeten we ongeveer testen?X"

# Remove the bookmark from its old location (end of the old last
# paragraph) so it can be re-added at the new location.
$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$newPara = $d.Paragraphs.Last
$sentinelStart = $newPara.Range.End - 2
$anchor = $d.Range($sentinelStart, $sentinelStart)
$d.Bookmarks.Add("_GoBack", $anchor)

# Strip the sentinel character back out again.
$sentinelRange = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinelRange.Text = ""
